$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting 2022/2023 down
$ws.Rows.Item(2).Insert()

# Clear the formatting that got carried into the new row so A2 keeps no style
$ws.Range("A2").ClearFormats()

# Set the new cell A2 to 2021
$ws.Range("A2").Value = 2021

# Select cell A2 as the active cell (reflected in sheetView selection)
$ws.Range("A2").Select()
